$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Rename "Photon Framework" -> "Photon Networking" in the project
# description line "... (C#, Unity Engine, Photon Framework)".
# ------------------------------------------------------------------

# Locate "Photon Framework" to get the precise character offsets of
# the trailing word "Framework" that needs to become "Networking".
$rngFull = $d.Content
$rngFull.Find.Execute("Photon Framework")
$frameworkStart = $rngFull.End - 9
$frameworkEnd = $rngFull.End

# The run immediately before it (" (C#") has identical character
# formatting, so normalisation would otherwise merge it into the
# edited run. "Freeze" it (toggle Bold off/on - a no-op formatting
# change) so it keeps its own run.
$rngPrefix = $d.Content
$rngPrefix.Find.Execute(" (C#")
$rngPrefix.Bold = $false
$rngPrefix.Bold = $true

# Replace the word itself.
$rngWord = $d.Range($frameworkStart, $frameworkEnd)
$rngWord.Text = "Networking"
$newEnd = $frameworkStart + 10

# Freeze the freshly written "Networking" run the same way, so it
# stays a distinct run instead of re-merging with its neighbours.
$rngNew = $d.Range($frameworkStart, $newEnd)
$rngNew.Bold = $false
$rngNew.Bold = $true

# Re-apply the freeze to the "(C#" run once more, since the edit
# above can re-trigger paragraph-wide run normalisation.
$rngPrefix2 = $d.Content
$rngPrefix2.Find.Execute(" (C#")
$rngPrefix2.Bold = $false
$rngPrefix2.Bold = $true

# ------------------------------------------------------------------
# Word keeps a single "_GoBack" bookmark marking the most recent
# edit location - move it from its old spot (after "'s annual") to
# right after the text we just typed ("Networking").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks.Item("_GoBack").Delete()
}
$rngBm = $d.Range($newEnd, $newEnd)
$d.Bookmarks.Add("_GoBack", $rngBm)
